$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Amazing Slime\UP\406 number A.mp4"'
$ws.Cells.Item(3, 2).Value = 'zzTESTzz'
$ws.Cells.Item(3, 3).Value = 'The king of the pirates'
$ws.Cells.Item(3, 4).Value = 'To make holding a key convenient, the hold() function can be used as a context manager and passed a string from the pyautogui.KEYBOARD_KEYS such as shift, ctrl, alt, and this key will be held for the duration of the with context block. See KEYBOARD_KEYS.'
$ws.Cells.Item(3, 5).Value = '15:31'
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = '11/9/2026'
$ws.Cells.Item(3, 6).ClearFormats()
$ws.Cells.Item(3, 7).Value = 'Uploaded'
$ws.Cells.Item(3, 8).Value = 'https://youtube.com/shorts/1q90a0XStHA'
